$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns, copying the header formatting
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record data for each player row (2-49)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 66  # AD
    $ws.Cells.Item($r, 31).Value = 96  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
